# Updated main GSC export data.
#
# The GSC Coverage export rolled forward by one day: the 2025-10-05 row
# (which had no data yet) is dropped from the "Chart" sheet, and every
# later date's row slides up to take its place - so each date keeps its
# own Not-indexed / Indexed / Impressions figures, the table just loses
# its leading (still-empty) day and shortens from 80 data rows to 79
# (A1:D81 -> A1:D80). The other sheets (Critical issues, Non-critical
# issues, Metadata) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 1 is the header (Date / Not indexed / Indexed / Impressions).
# Row 2 is 2025-10-05 with no data yet (blank Not-indexed/Indexed, 0
# impressions) - delete it and let every following row shift up one.
$ws.Rows.Item(2).Delete()
